$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.05052846209076733
$ws.Range("D2").Value = 0.008943546438445082
$ws.Range("E2").Value = 0.4304237066138796
$ws.Range("F2").Value = 1.038957906139956
$ws.Range("G2").Value = 0.002342237695844014
$ws.Range("N2").Value = 2.873581687696344
$ws.Range("O2").Value = 3.581084263747925

$ws.Range("C3").Value = 0.04476206481484724
$ws.Range("D3").Value = 0.007995848526245908
$ws.Range("E3").Value = 0.3748685171766653
$ws.Range("F3").Value = 0.9555018548654601
$ws.Range("G3").Value = 0.00234811694982799
$ws.Range("N3").Value = 2.562605684679454
$ws.Range("O3").Value = 3.290152896055872

$ws.Range("C4").Value = 0.04124261813031183
$ws.Range("D4").Value = 0.007420112891711739
$ws.Range("E4").Value = 0.3409164150051254
$ws.Range("F4").Value = 0.9049358929064368
$ws.Range("G4").Value = 0.002351912108017545
$ws.Range("N4").Value = 2.371325805375761
$ws.Range("O4").Value = 3.113860341885299

$ws.Range("C5").Value = 0.03981353979416724
$ws.Range("D5").Value = 0.007186965551756685
$ws.Range("E5").Value = 0.3271169890140158
$ws.Range("F5").Value = 0.8844966276862749
$ws.Range("G5").Value = 0.002353505441391185
$ws.Range("N5").Value = 2.293303068607429
$ws.Range("O5").Value = 3.042596276585243

$ws.Range("C6").Value = 0.03957654638250574
$ws.Range("D6").Value = 0.007148338117151098
$ws.Range("E6").Value = 0.3248276933736634
$ws.Range("F6").Value = 0.8811126802679041
$ws.Range("G6").Value = 0.002353772843655402
$ws.Range("N6").Value = 2.280343261403573
$ws.Range("O6").Value = 3.030797410986395

$ws.Range("C7").Value = 0.04122332458327094
$ws.Range("D7").Value = 0.007416962736222388
$ws.Range("E7").Value = 0.3407301688095856
$ws.Range("F7").Value = 0.9046595708997387
$ws.Range("G7").Value = 0.002351933406542927
$ws.Range("N7").Value = 2.370273851392824
$ws.Range("O7").Value = 3.112896930884858

$ws.Range("C8").Value = 0.04853570168489796
$ws.Range("D8").Value = 0.008615455172794384
$ws.Range("E8").Value = 0.4112330468379781
$ws.Range("F8").Value = 1.010039843139111
$ws.Range("G8").Value = 0.002344226522347518
$ws.Range("N8").Value = 2.766433886209825
$ws.Range("O8").Value = 3.480278027715315

$ws.Range("C9").Value = 0.06305247523309276
$ws.Range("D9").Value = 0.01101820639424744
$ws.Range("E9").Value = 0.5509232592974627
$ws.Range("F9").Value = 1.222218960642749
$ws.Range("G9").Value = 0.002330574835277527
$ws.Range("N9").Value = 3.540180268007646
$ws.Range("O9").Value = 4.219864929822393

$ws.Range("C10").Value = 0.07384021095181481
$ws.Range("D10").Value = 0.0128208404012895
$ws.Range("E10").Value = 0.6546786155740136
$ws.Range("F10").Value = 1.381719871913816
$ws.Range("G10").Value = 0.002321423881673422
$ws.Range("N10").Value = 4.10623028343673
$ws.Range("O10").Value = 4.775787231787206

$ws.Range("C11").Value = 0.07877756143588499
$ws.Range("D11").Value = 0.01365020014181084
$ws.Range("E11").Value = 0.7021797081820296
$ws.Range("F11").Value = 1.455119823126978
$ws.Range("G11").Value = 0.002317449146422014
$ws.Range("N11").Value = 4.363110593465422
$ws.Range("O11").Value = 5.03161050597248

$ws.Range("C12").Value = 0.08065175668153302
$ws.Range("D12").Value = 0.01396569544306914
$ws.Range("E12").Value = 0.7202150718742786
$ws.Range("F12").Value = 1.48303976134892
$ws.Range("G12").Value = 0.002315970862815239
$ws.Range("N12").Value = 4.460285735714251
$ws.Range("O12").Value = 5.128920722311477

$ws.Range("C13").Value = 0.08024791027808931
$ws.Range("D13").Value = 0.01389768257172364
$ws.Range("E13").Value = 0.71632864187238
$ws.Range("F13").Value = 1.477021082648264
$ws.Range("G13").Value = 0.002316288045991044
$ws.Range("N13").Value = 4.439361943450422
$ws.Range("O13").Value = 5.107943632011711

$ws.Range("C14").Value = 0.07893166058704537
$ws.Range("D14").Value = 0.01367612683555564
$ws.Range("E14").Value = 0.7036625070413294
$ws.Range("F14").Value = 1.45741428624936
$ws.Range("G14").Value = 0.002317326989631076
$ws.Range("N14").Value = 4.371107314139522
$ws.Range("O14").Value = 5.039607467092196

$ws.Range("C15").Value = 0.07812601590303814
$ws.Range("D15").Value = 0.01354060719010874
$ws.Range("E15").Value = 0.6959104778898961
$ws.Range("F15").Value = 1.445420948117004
$ws.Range("G15").Value = 0.002317966865639248
$ws.Range("N15").Value = 4.329286057409945
$ws.Range("O15").Value = 4.997806726364843

$ws.Range("C16").Value = 0.07351816499694053
$ws.Range("D16").Value = 0.0127668354652215
$ws.Range("E16").Value = 0.6515807523204415
$ws.Range("F16").Value = 1.376940287943825
$ws.Range("G16").Value = 0.002321687408931128
$ws.Range("N16").Value = 4.089429168003562
$ws.Range("O16").Value = 4.759128753694995

$ws.Range("C17").Value = 0.07069922541701601
$ws.Range("D17").Value = 0.0122946071648613
$ws.Range("E17").Value = 0.6244663244897453
$ws.Range("F17").Value = 1.335148049666174
$ws.Range("G17").Value = 0.002324017883012411
$ws.Range("N17").Value = 3.94211849063862
$ws.Range("O17").Value = 4.613468198980399

$ws.Range("C18").Value = 0.06908065205502112
$ws.Range("D18").Value = 0.01202386667441147
$ws.Range("E18").Value = 0.6088990123416664
$ws.Range("F18").Value = 1.311189294299879
$ws.Range("G18").Value = 0.002325376024265097
$ws.Range("N18").Value = 3.857331695637754
$ws.Range("O18").Value = 4.529963168608333

$ws.Range("C19").Value = 0.06853310623974096
$ws.Range("D19").Value = 0.01193234570065016
$ws.Range("E19").Value = 0.6036329010447616
$ws.Range("F19").Value = 1.303090729238761
$ws.Range("G19").Value = 0.002325838915274946
$ws.Range("N19").Value = 3.828614786364199
$ws.Range("O19").Value = 4.501736622715214

$ws.Range("C20").Value = 0.07099901434185085
$ws.Range("D20").Value = 0.01234478568275676
$ws.Range("E20").Value = 0.6273497534900514
$ws.Range("F20").Value = 1.339588697244238
$ws.Range("G20").Value = 0.00232376796816338
$ws.Range("N20").Value = 3.95780600327754
$ws.Range("O20").Value = 4.628945444994827

$ws.Range("C21").Value = 0.07931815037034085
$ws.Range("D21").Value = 0.01374116343618681
$ws.Range("E21").Value = 0.7073815280194395
$ws.Range("F21").Value = 1.463169854790635
$ws.Range("G21").Value = 0.002317021099257929
$ws.Range("N21").Value = 4.391158149571083
$ws.Range("O21").Value = 5.059667523342455

$ws.Range("C22").Value = 0.08478168316780454
$ws.Range("D22").Value = 0.01466218457728274
$ws.Range("E22").Value = 0.7599671772599095
$ws.Range("F22").Value = 1.544667555240125
$ws.Range("G22").Value = 0.002312768120562306
$ws.Range("N22").Value = 4.67379181795809
$ws.Range("O22").Value = 5.343714382450571

$ws.Range("C23").Value = 0.08186319461719904
$ws.Range("D23").Value = 0.01416981774028869
$ws.Range("E23").Value = 0.7318741492869663
$ws.Range("F23").Value = 1.501102581173711
$ws.Range("G23").Value = 0.002315023756898066
$ws.Range("N23").Value = 4.523002190001307
$ws.Range("O23").Value = 5.19187566491928

$ws.Range("C24").Value = 0.07086347333060417
$ws.Range("D24").Value = 0.01232209764616954
$ws.Range("E24").Value = 0.6260460897340607
$ws.Range("F24").Value = 1.33758086846484
$ws.Range("G24").Value = 0.002323880897580954
$ws.Range("N24").Value = 3.950713976768498
$ws.Range("O24").Value = 4.621947444414047

$ws.Range("C25").Value = 0.05910487357310501
$ws.Range("D25").Value = 0.01036203391920765
$ws.Range("E25").Value = 0.5129541931939059
$ws.Range("F25").Value = 1.164202196230377
$ws.Range("G25").Value = 0.002334112758888872
$ws.Range("N25").Value = 3.331249627311138
$ws.Range("O25").Value = 4.017647422946709
